# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# 1) "ODI Batting" sheet: two leftover empty INNING_NUMBER cells (B5 / B10)
#    are cleared out (they were already blank, this just drops the dead
#    cell record).
# 2) A brand-new "ODI Batting Extra" sheet is appended after "ODI Bowling"
#    with the scraped extra-batting-stat columns.

$wb = $excel.ActiveWorkbook

# --- 1. ODI Batting: drop the two stray empty B cells -----------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B5").ClearContents()
$batting.Range("B10").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet after "ODI Bowling" -----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Reuse the bold/bordered/centered header style (style index 1) that's
# already used for header rows on the other sheets.
$headerStyleSource = $wb.Worksheets.Item("Player Info")
$headerStyleSource.Range("A1:D1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Row 2 - match 4258 (no batting-position/num4/num6/percent data scraped)
$ws.Range("A2").Value = "'4258"
$ws.Range("F2").Value = "NO"

# Row 3 - match 4268
$ws.Range("A3").Value = "'4268"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = "'5"
$ws.Range("D3").Value = "'6"
$ws.Range("E3").Value = "'23.40%"
$ws.Range("F3").Value = "YES"

# Row 4 - match 4270
$ws.Range("A4").Value = "'4270"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = "'2"
$ws.Range("D4").Value = "'1"
$ws.Range("E4").Value = "'7.35%"
$ws.Range("F4").Value = "NO"

# Row 5 - match 4398 (no num4/num6/percent data scraped)
$ws.Range("A5").Value = "'4398"
$ws.Range("B5").Value = 6
$ws.Range("F5").Value = "NO"

# Row 6 - match 4399
$ws.Range("A6").Value = "'4399"
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = "'0"
$ws.Range("D6").Value = "'0"
$ws.Range("E6").Value = "'4.28%"
$ws.Range("F6").Value = "NO"

# Row 7 - match 4400
$ws.Range("A7").Value = "'4400"
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = "'0"
$ws.Range("D7").Value = "'0"
$ws.Range("E7").Value = "'1.40%"
$ws.Range("F7").Value = "NO"

# Row 8 - match 4483
$ws.Range("A8").Value = "'4483"
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = "'2"
$ws.Range("D8").Value = "'2"
$ws.Range("E8").Value = "'19.44%"
$ws.Range("F8").Value = "NO"

# Row 9 - match 4484
$ws.Range("A9").Value = "'4484"
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = "'0"
$ws.Range("D9").Value = "'0"
$ws.Range("E9").Value = "'0.53%"
$ws.Range("F9").Value = "NO"

# Row 10 - match 4486 (no batting-position/num4/num6/percent data scraped)
$ws.Range("A10").Value = "'4486"
$ws.Range("F10").Value = "NO"

# Restore the original active sheet (adding a sheet shifts selection to it)
$wb.Worksheets.Item("Player Info").Activate()
